$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 287-288; everything currently at row 287
# onward (through the old last row 406) shifts down to 289-408.
$ws.Rows("287:288").Insert()

# Populate the two newly inserted rows with the new weekly records.
$ws.Range("A287").Value = 5
$ws.Range("B287").Value = "Macroferia Regional de Talca"
$ws.Range("C287").Value = "Maule"
$ws.Range("D287").Value = 44900
$ws.Range("E287").Value = 7
$ws.Range("F287").Value = 100112006
$ws.Range("G287").Value = "Repollo"
$ws.Range("H287").Value = "Copenhague"
$ws.Range("I287").Value = "Segunda"
$ws.Range("J287").Value = 3000
$ws.Range("K287").Value = 700
$ws.Range("L287").Value = 700
$ws.Range("M287").Value = 700
$ws.Range("N287").Value = "`$/unidad"
$ws.Range("O287").Value = "Provincia del Elquí"
$ws.Range("P287").Value = 700
$ws.Range("Q287").Value = 1
$ws.Range("R287").Value = "Hortaliza"

$ws.Range("A288").Value = 5
$ws.Range("B288").Value = "Macroferia Regional de Talca"
$ws.Range("C288").Value = "Maule"
$ws.Range("D288").Value = 44900
$ws.Range("E288").Value = 7
$ws.Range("F288").Value = 100112006
$ws.Range("G288").Value = "Repollo"
$ws.Range("H288").Value = "Crespo record"
$ws.Range("I288").Value = "Segunda"
$ws.Range("J288").Value = 4000
$ws.Range("K288").Value = 700
$ws.Range("L288").Value = 700
$ws.Range("M288").Value = 700
$ws.Range("N288").Value = "`$/unidad"
$ws.Range("O288").Value = "Provincia del Elquí"
$ws.Range("P288").Value = 700
$ws.Range("Q288").Value = 1
$ws.Range("R288").Value = "Hortaliza"
